$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Efna1"
$ws.Cells.Item(2,3).Value = "Epha3"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 23.73148533333334
$ws.Cells.Item(2,8).Value = 71.194456
$ws.Cells.Item(2,9).Value = 0.8653076146801144
$ws.Cells.Item(2,10).Value = 0.8653076146801145
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.07579599999999999
$ws.Cells.Item(2,14).Value = 0.227388
$ws.Cells.Item(2,15).Value = 0.001780200955210419
$ws.Cells.Item(2,16).Value = 0.001780200955210419
$ws.Cells.Item(2,17).Value = 1.798751662325333
$ws.Cells.Item(2,18).Value = 16.188764960928
$ws.Cells.Item(2,19).Value = 0.001540421442204388
$ws.Cells.Item(2,20).Value = 0.001540421442204389

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Efna1"
$ws.Cells.Item(3,3).Value = "Epha3"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 23.73148533333334
$ws.Cells.Item(3,8).Value = 71.194456
$ws.Cells.Item(3,9).Value = 0.8653076146801144
$ws.Cells.Item(3,10).Value = 0.8653076146801145
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 40.12734033333333
$ws.Cells.Item(3,14).Value = 120.382021
$ws.Cells.Item(3,15).Value = 0.9424604146848589
$ws.Cells.Item(3,16).Value = 0.9424604146848587
$ws.Cells.Item(3,17).Value = 952.2813885861751
$ws.Cells.Item(3,18).Value = 8570.532497275575
$ws.Cells.Item(3,19).Value = 0.8155181733613867
$ws.Cells.Item(3,20).Value = 0.8155181733613867

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Efna1"
$ws.Cells.Item(4,3).Value = "Epha3"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 23.73148533333334
$ws.Cells.Item(4,8).Value = 71.194456
$ws.Cells.Item(4,9).Value = 0.8653076146801144
$ws.Cells.Item(4,10).Value = 0.8653076146801145
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.344072666666666
$ws.Cells.Item(4,14).Value = 7.032217999999999
$ws.Cells.Item(4,15).Value = 0.05505462557763778
$ws.Cells.Item(4,16).Value = 0.05505462557763778
$ws.Cells.Item(4,17).Value = 55.62832610926755
$ws.Cells.Item(4,18).Value = 500.654934983408
$ws.Cells.Item(4,19).Value = 0.04763918673569256
$ws.Cells.Item(4,20).Value = 0.04763918673569257

$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Efna1"
$ws.Cells.Item(5,3).Value = "Epha3"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 23.73148533333334
$ws.Cells.Item(5,8).Value = 71.194456
$ws.Cells.Item(5,9).Value = 0.8653076146801144
$ws.Cells.Item(5,10).Value = 0.8653076146801145
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.03000666666666667
$ws.Cells.Item(5,14).Value = 0.09002
$ws.Cells.Item(5,15).Value = 0.0007047587822930054
$ws.Cells.Item(5,16).Value = 0.0007047587822930053
$ws.Cells.Item(5,17).Value = 0.7121027699022223
$ws.Cells.Item(5,18).Value = 6.40892492912
$ws.Cells.Item(5,19).Value = 0.0006098331408308225
$ws.Cells.Item(5,20).Value = 0.0006098331408308225

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Efna1"
$ws.Cells.Item(6,3).Value = "Epha3"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.471212
$ws.Cells.Item(6,8).Value = 7.413636
$ws.Cells.Item(6,9).Value = 0.09010639372350319
$ws.Cells.Item(6,10).Value = 0.09010639372350321
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.07579599999999999
$ws.Cells.Item(6,14).Value = 0.227388
$ws.Cells.Item(6,15).Value = 0.001780200955210419
$ws.Cells.Item(6,16).Value = 0.001780200955210419
$ws.Cells.Item(6,17).Value = 0.187307984752
$ws.Cells.Item(6,18).Value = 1.685771862768
$ws.Cells.Item(6,19).Value = 0.0001604074881771465
$ws.Cells.Item(6,20).Value = 0.0001604074881771465

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Efna1"
$ws.Cells.Item(7,3).Value = "Epha3"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.471212
$ws.Cells.Item(7,8).Value = 7.413636
$ws.Cells.Item(7,9).Value = 0.09010639372350319
$ws.Cells.Item(7,10).Value = 0.09010639372350321
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 40.12734033333333
$ws.Cells.Item(7,14).Value = 120.382021
$ws.Cells.Item(7,15).Value = 0.9424604146848589
$ws.Cells.Item(7,16).Value = 0.9424604146848587
$ws.Cells.Item(7,17).Value = 99.16316495981732
$ws.Cells.Item(7,18).Value = 892.468484638356
$ws.Cells.Item(7,19).Value = 0.08492170919440999
$ws.Cells.Item(7,20).Value = 0.08492170919440999

$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Efna1"
$ws.Cells.Item(8,3).Value = "Epha3"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 2.471212
$ws.Cells.Item(8,8).Value = 7.413636
$ws.Cells.Item(8,9).Value = 0.09010639372350319
$ws.Cells.Item(8,10).Value = 0.09010639372350321
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 2.344072666666666
$ws.Cells.Item(8,14).Value = 7.032217999999999
$ws.Cells.Item(8,15).Value = 0.05505462557763778
$ws.Cells.Item(8,16).Value = 0.05505462557763778
$ws.Cells.Item(8,17).Value = 5.792700502738666
$ws.Cells.Item(8,18).Value = 52.13430452464799
$ws.Cells.Item(8,19).Value = 0.004960773768598679
$ws.Cells.Item(8,20).Value = 0.00496077376859868

$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Efna1"
$ws.Cells.Item(9,3).Value = "Epha3"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 2.471212
$ws.Cells.Item(9,8).Value = 7.413636
$ws.Cells.Item(9,9).Value = 0.09010639372350319
$ws.Cells.Item(9,10).Value = 0.09010639372350321
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.03000666666666667
$ws.Cells.Item(9,14).Value = 0.09002
$ws.Cells.Item(9,15).Value = 0.0007047587822930054
$ws.Cells.Item(9,16).Value = 0.0007047587822930053
$ws.Cells.Item(9,17).Value = 0.07415283474666667
$ws.Cells.Item(9,18).Value = 0.6673755127200001
$ws.Cells.Item(9,19).Value = [double]"6.350327231739021E-05"
$ws.Cells.Item(9,20).Value = [double]"6.350327231739021E-05"

$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Efna1"
$ws.Cells.Item(10,3).Value = "Epha3"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.148663
$ws.Cells.Item(10,8).Value = 3.445989
$ws.Cells.Item(10,9).Value = 0.04188304383987305
$ws.Cells.Item(10,10).Value = 0.04188304383987305
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.07579599999999999
$ws.Cells.Item(10,14).Value = 0.227388
$ws.Cells.Item(10,15).Value = 0.001780200955210419
$ws.Cells.Item(10,16).Value = 0.001780200955210419
$ws.Cells.Item(10,17).Value = 0.08706406074799998
$ws.Cells.Item(10,18).Value = 0.7835765467319999
$ws.Cells.Item(10,19).Value = [double]"7.456023465086184E-05"
$ws.Cells.Item(10,20).Value = [double]"7.456023465086184E-05"

$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Efna1"
$ws.Cells.Item(11,3).Value = "Epha3"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 1.148663
$ws.Cells.Item(11,8).Value = 3.445989
$ws.Cells.Item(11,9).Value = 0.04188304383987305
$ws.Cells.Item(11,10).Value = 0.04188304383987305
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 40.12734033333333
$ws.Cells.Item(11,14).Value = 120.382021
$ws.Cells.Item(11,15).Value = 0.9424604146848589
$ws.Cells.Item(11,16).Value = 0.9424604146848587
$ws.Cells.Item(11,17).Value = 46.09279112930766
$ws.Cells.Item(11,18).Value = 414.835120163769
$ws.Cells.Item(11,19).Value = 0.03947311086559088
$ws.Cells.Item(11,20).Value = 0.03947311086559087

$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Efna1"
$ws.Cells.Item(12,3).Value = "Epha3"
$ws.Cells.Item(12,4).Value = "MuSCs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 1.148663
$ws.Cells.Item(12,8).Value = 3.445989
$ws.Cells.Item(12,9).Value = 0.04188304383987305
$ws.Cells.Item(12,10).Value = 0.04188304383987305
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 2.344072666666666
$ws.Cells.Item(12,14).Value = 7.032217999999999
$ws.Cells.Item(12,15).Value = 0.05505462557763778
$ws.Cells.Item(12,16).Value = 0.05505462557763778
$ws.Cells.Item(12,17).Value = 2.692549541511333
$ws.Cells.Item(12,18).Value = 24.232945873602
$ws.Cells.Item(12,19).Value = 0.002305855296655999
$ws.Cells.Item(12,20).Value = 0.002305855296655999

$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Efna1"
$ws.Cells.Item(13,3).Value = "Epha3"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 1.148663
$ws.Cells.Item(13,8).Value = 3.445989
$ws.Cells.Item(13,9).Value = 0.04188304383987305
$ws.Cells.Item(13,10).Value = 0.04188304383987305
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 0.6666666666666666
$ws.Cells.Item(13,13).Value = 0.03000666666666667
$ws.Cells.Item(13,14).Value = 0.09002
$ws.Cells.Item(13,15).Value = 0.0007047587822930054
$ws.Cells.Item(13,16).Value = 0.0007047587822930053
$ws.Cells.Item(13,17).Value = 0.03446754775333333
$ws.Cells.Item(13,18).Value = 0.31020792978
$ws.Cells.Item(13,19).Value = [double]"2.951744297531349E-05"
$ws.Cells.Item(13,20).Value = [double]"2.951744297531349E-05"

$ws.Cells.Item(14,1).Value = "Resolving-Mac"
$ws.Cells.Item(14,2).Value = "Efna1"
$ws.Cells.Item(14,3).Value = "Epha3"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 1
$ws.Cells.Item(14,6).Value = 0.3333333333333333
$ws.Cells.Item(14,7).Value = 0.07412966666666666
$ws.Cells.Item(14,8).Value = 0.222389
$ws.Cells.Item(14,9).Value = 0.002702947756509242
$ws.Cells.Item(14,10).Value = 0.002702947756509243
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 0.6666666666666666
$ws.Cells.Item(14,13).Value = 0.07579599999999999
$ws.Cells.Item(14,14).Value = 0.227388
$ws.Cells.Item(14,15).Value = 0.001780200955210419
$ws.Cells.Item(14,16).Value = 0.001780200955210419
$ws.Cells.Item(14,17).Value = 0.005618732214666666
$ws.Cells.Item(14,18).Value = 0.050568589932
$ws.Cells.Item(14,19).Value = [double]"4.811790178021611E-06"
$ws.Cells.Item(14,20).Value = [double]"4.811790178021612E-06"

$ws.Cells.Item(15,1).Value = "Resolving-Mac"
$ws.Cells.Item(15,2).Value = "Efna1"
$ws.Cells.Item(15,3).Value = "Epha3"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 1
$ws.Cells.Item(15,6).Value = 0.3333333333333333
$ws.Cells.Item(15,7).Value = 0.07412966666666666
$ws.Cells.Item(15,8).Value = 0.222389
$ws.Cells.Item(15,9).Value = 0.002702947756509242
$ws.Cells.Item(15,10).Value = 0.002702947756509243
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 40.12734033333333
$ws.Cells.Item(15,14).Value = 120.382021
$ws.Cells.Item(15,15).Value = 0.9424604146848589
$ws.Cells.Item(15,16).Value = 0.9424604146848587
$ws.Cells.Item(15,17).Value = 2.974626363129889
$ws.Cells.Item(15,18).Value = 26.771637268169
$ws.Cells.Item(15,19).Value = 0.002547421263471209
$ws.Cells.Item(15,20).Value = 0.002547421263471209

$ws.Cells.Item(16,1).Value = "Resolving-Mac"
$ws.Cells.Item(16,2).Value = "Efna1"
$ws.Cells.Item(16,3).Value = "Epha3"
$ws.Cells.Item(16,4).Value = "MuSCs"
$ws.Cells.Item(16,5).Value = 1
$ws.Cells.Item(16,6).Value = 0.3333333333333333
$ws.Cells.Item(16,7).Value = 0.07412966666666666
$ws.Cells.Item(16,8).Value = 0.222389
$ws.Cells.Item(16,9).Value = 0.002702947756509242
$ws.Cells.Item(16,10).Value = 0.002702947756509243
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 2.344072666666666
$ws.Cells.Item(16,14).Value = 7.032217999999999
$ws.Cells.Item(16,15).Value = 0.05505462557763778
$ws.Cells.Item(16,16).Value = 0.05505462557763778
$ws.Cells.Item(16,17).Value = 0.1737653254224444
$ws.Cells.Item(16,18).Value = 1.563887928802
$ws.Cells.Item(16,19).Value = 0.0001488097766905324
$ws.Cells.Item(16,20).Value = 0.0001488097766905324

$ws.Cells.Item(17,1).Value = "Resolving-Mac"
$ws.Cells.Item(17,2).Value = "Efna1"
$ws.Cells.Item(17,3).Value = "Epha3"
$ws.Cells.Item(17,4).Value = "Resolving-Mac"
$ws.Cells.Item(17,5).Value = 1
$ws.Cells.Item(17,6).Value = 0.3333333333333333
$ws.Cells.Item(17,7).Value = 0.07412966666666666
$ws.Cells.Item(17,8).Value = 0.222389
$ws.Cells.Item(17,9).Value = 0.002702947756509242
$ws.Cells.Item(17,10).Value = 0.002702947756509243
$ws.Cells.Item(17,11).Value = 2
$ws.Cells.Item(17,12).Value = 0.6666666666666666
$ws.Cells.Item(17,13).Value = 0.03000666666666667
$ws.Cells.Item(17,14).Value = 0.09002
$ws.Cells.Item(17,15).Value = 0.0007047587822930054
$ws.Cells.Item(17,16).Value = 0.0007047587822930053
$ws.Cells.Item(17,17).Value = 0.002224384197777778
$ws.Cells.Item(17,18).Value = 0.02001945778
$ws.Cells.Item(17,19).Value = [double]"1.904926169479064E-06"
$ws.Cells.Item(17,20).Value = [double]"1.904926169479065E-06"

